$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 13 (which only held the "Docentes responsáveis:" value in B/C,
# with no label in column A) - this shifts rows 14-22 up to 13-21, carrying
# their row heights/styles along.
$ws.Rows("13").Delete()

# Targeted content fixes on top of the shift
$ws.Range("B10").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C10").Value = "5840671 - Francisco José Moreira Chaves"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Formula = "'01/01/2012"
$ws.Range("C15").Formula = "'01/01/2012"

$ws.Range("B18").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C18").Value = "5840671 - Francisco José Moreira Chaves"

$ws.Range("B19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."

$ws.Range("B20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("C20").Value = "MF = (0,30*P1 + 0,60*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
